# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.7127328510149897; E = 0.4998867070740569; F = 0; G = 6.048734245549538 }
    3 = @{ B = 0.7287194209349384; C = 0.3375848360084654;  D = 0.7127328510149897; E = 0.4998867070740569; F = 0; G = 2.27892381503245 }
    4 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.1529057820181812; E = 0.4998867070740569; F = 1; G = 5.488907176552729 }
    5 = @{ B = 0.3464964993005633; C = 0.3375848360084654;  D = 0.1529057820181812; E = 0.4998867070740569; F = 0; G = 1.336873824401267 }
    6 = @{ B = 3.182878228561681;  C = 1.65323645889881;   D = 0.1529057820181812; E = 0.4998867070740569; F = 1; G = 5.488907176552729 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}
